# "get rid of the test data" — clear the sample/test rows (2-6) from the
# client_persons mapping sheet, leaving only the header row (row 1) intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the block of test data (rows 2-6, all used columns) and clear it,
# the way a user would select the rows in the UI and hit Delete.
$range = $ws.Range("A2:G6")
$range.Select()
$range.ClearContents()
